$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# D5: new "1. Back Bar" label next to the existing "Back Bar" row (A5)
$ws.Range("D5").Value = "1. Back Bar"

# D6: new detailed label next to the existing "Menu" row (A6)
$enDash = [char]0x2013
$ws.Range("D6").Value = "1. Main Menu, Menu " + $enDash + " Features, 3. Spirits List Menu, 4. Bottle List Menu, 5. After Dinner / Dessert Menu"

# Column A and D widen slightly to fit the new content
$ws.Columns("A").ColumnWidth = 20.333333333333332
$ws.Columns("D").ColumnWidth = 117.0

# Active cell/selection moves to D11
$ws.Range("D11").Select()
